# Updated symbol list on Sat Dec 24 10:21:19 UTC 2022 with GitHub Actions
#
# Applies the price/coin-list refresh to the "cryptos" worksheet.
# Columns: A=idx, B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must be stored as literal text (not
# re-interpreted/renormalized as a floating point number), while leaving the
# cell's original style/number-format untouched.
function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $c = $ws.Range($CellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $Value
    $c.Style = $origStyle
}

# Helper: update an entire data row (Coin, Link, Price, Volume label)
function Set-Row {
    param(
        [int]$Row,
        [string]$Coin,
        [string]$Link,
        [string]$Price,
        [string]$Volume
    )
    $ws.Range("B$Row").Value = $Coin
    $ws.Range("C$Row").Value = $Link
    Set-TextValue "D$Row" $Price
    $ws.Range("E$Row").Value = $Volume
}

# --- Simple price-only updates -------------------------------------------
Set-TextValue "D2" "244.85"
Set-TextValue "D3" "21.98"
Set-TextValue "D4" "5.392"
Set-TextValue "D6" "3.389"

# --- Rows 7-23: coin ranking reshuffled (each row's data shifts up one,
#     with a refreshed KuCoinToken entry appended at the former LEO slot) ---
Set-Row 7  "MXToken"                             "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"                                 "0.8115"    "6MXTokenMX"
Set-Row 8  "FTXToken"                            "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"                                  "0.9549"    "7FTXTokenFTT"
Set-Row 9  "WazirX"                              "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"                                     "0.1428"    "8WazirXWRX"
Set-Row 10 "MandalaExchangeToken"                "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"                  "0.07405"   "9MandalaExchangeTokenMDX"
Set-Row 11 "LiechtensteinCryptoassetsExchange"   "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"          "0.03381"   "10LiechtensteinCryptoassetsExchangeLCX"
Set-Row 12 "BitrueCoin"                          "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"                                "0.03055"   "11BitrueCoinBTR"
Set-Row 13 "BitMartToken"                        "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"                          "0.09426"   "12BitMartTokenBMX"
Set-Row 14 "MCDex"                               "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"                                     "4.002"     "13MCDexMCB"
Set-Row 15 "BitForexToken"                       "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"                          "0.001596"  "14BitForexTokenBF"
Set-Row 16 "CoinExToken"                         "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"                           "0.04795"   "15CoinExTokenCET"
Set-Row 17 "One"                                 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"                                    "0.0005872" "16OneONE"
Set-Row 18 "TigerCash"                           "https://coinranking.com/coin/6hIn06L2+tigercash-tch"                                  "0.006165"  "17TigerCashTCH"
Set-Row 19 "HotbitToken"                         "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"                            "0.005069"  "18HotbitTokenHTB"
Set-Row 20 "BitKan"                              "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"                                "0.0009870" "19BitKanKAN"
Set-Row 21 "NitroEx"                             "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"                                 "0.0001000" "20NitroExNTX"
Set-Row 22 "LEO"                                 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"                                    "3.702"     "21LEOLEO"
Set-Row 23 "KuCoinToken"                         "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"                           "6.394"     "22KuCoinTokenKCS"

# --- More simple price-only updates ---------------------------------------
Set-TextValue "D24" "2.186"
Set-TextValue "D25" "0.3231"
Set-TextValue "D40" "0.03991"

# --- Rows 41-43: another small reshuffle ----------------------------------
Set-Row 41 "BKEXToken"  "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"      "0.1075"    "40BKEXTokenBKK"
Set-Row 42 "CEJI"       "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"          "0.002721"  "41CEJICEJI"
Set-Row 43 "KickToken"  "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick" "0.003020"  "42KickTokenKICK"

# --- Remaining simple price updates ---------------------------------------
Set-TextValue "D44" "0.005867"
Set-TextValue "D45" "0.00005279"
Set-TextValue "D47" "0.8013"

# --- Row 48: price update + volume label change ---------------------------
Set-TextValue "D48" "0.02333"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
